$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1725239616613418
$ws.Range("C2").Value = 0.5846645367412141
$ws.Range("J2").Value = 0.01277955271565495
$ws.Range("P2").Value = 0.121405750798722
$ws.Range("S2").Value = 0.1086261980830671
$ws.Range("B3").Value = 0.01063829787234043
$ws.Range("C3").Value = 0.02659574468085106
$ws.Range("J3").Value = 0.0425531914893617
$ws.Range("P3").Value = 0.7393617021276596
$ws.Range("S3").Value = 0.1808510638297872
$ws.Range("J4").Value = 0.125
$ws.Range("P4").Value = 0.8
$ws.Range("S4").Value = 0.075
$ws.Range("B6").Value = 0.06274509803921569
$ws.Range("D6").Value = 0.00392156862745098
$ws.Range("F6").Value = 0.09803921568627451
$ws.Range("J6").Value = 0.2549019607843137
$ws.Range("O6").Value = 0.03137254901960784
$ws.Range("Q6").Value = 0.09803921568627451
$ws.Range("R6").Value = 0.05882352941176471
$ws.Range("S6").Value = 0.392156862745098
$ws.Range("B7").Value = 0.07246376811594203
$ws.Range("D7").Value = 0.01449275362318841
$ws.Range("F7").Value = 0.06280193236714976
$ws.Range("J7").Value = 0.1159420289855072
$ws.Range("O7").Value = 0.01932367149758454
$ws.Range("Q7").Value = 0.2270531400966184
$ws.Range("R7").Value = 0.09178743961352658
$ws.Range("S7").Value = 0.3961352657004831
$ws.Range("B8").Value = 0.108433734939759
$ws.Range("D8").Value = 0.01807228915662651
$ws.Range("E8").Value = 0.004016064257028112
$ws.Range("F8").Value = 0.07429718875502007
$ws.Range("J8").Value = 0.1044176706827309
$ws.Range("O8").Value = 0.01405622489959839
$ws.Range("Q8").Value = 0.1807228915662651
$ws.Range("R8").Value = 0.07028112449799197
$ws.Range("S8").Value = 0.4257028112449799
$ws.Range("B9").Value = 0.07537688442211055
$ws.Range("D9").Value = 0.02512562814070352
$ws.Range("F9").Value = 0.05527638190954774
$ws.Range("J9").Value = 0.1557788944723618
$ws.Range("O9").Value = 0.03015075376884422
$ws.Range("Q9").Value = 0.1959798994974874
$ws.Range("R9").Value = 0.09045226130653267
$ws.Range("S9").Value = 0.3718592964824121
$ws.Range("B10").Value = 0.1118568232662192
$ws.Range("D10").Value = 0.01789709172259508
$ws.Range("E10").Value = 0.001491424310216256
$ws.Range("F10").Value = 0.06935123042505593
$ws.Range("J10").Value = 0.1379567486950037
$ws.Range("O10").Value = 0.01789709172259508
$ws.Range("Q10").Value = 0.2132736763609247
$ws.Range("R10").Value = 0.07531692766592095
$ws.Range("S10").Value = 0.3549589858314691
$ws.Range("G11").Value = 0.1557632398753894
$ws.Range("J11").Value = 0.07476635514018691
$ws.Range("K11").Value = 0.205607476635514
$ws.Range("L11").Value = 0.5451713395638629
$ws.Range("S11").Value = 0.01869158878504673
$ws.Range("G12").Value = 0.7336956521739131
$ws.Range("J12").Value = 0.1793478260869565
$ws.Range("K12").Value = 0.02173913043478261
$ws.Range("L12").Value = 0.04891304347826087
$ws.Range("S12").Value = 0.01630434782608696
$ws.Range("G13").Value = 0.65
$ws.Range("J13").Value = 0.25
$ws.Range("S13").Value = 0.1
$ws.Range("F15").Value = 0.0205761316872428
$ws.Range("H15").Value = 0.1234567901234568
$ws.Range("I15").Value = 0.06995884773662552
$ws.Range("J15").Value = 0.3539094650205761
$ws.Range("K15").Value = 0.06995884773662552
$ws.Range("M15").Value = 0.00823045267489712
$ws.Range("O15").Value = 0.102880658436214
$ws.Range("S15").Value = 0.2510288065843622
$ws.Range("F16").Value = 0.0198019801980198
$ws.Range("H16").Value = 0.2128712871287129
$ws.Range("I16").Value = 0.04455445544554455
$ws.Range("J16").Value = 0.3960396039603961
$ws.Range("K16").Value = 0.1138613861386139
$ws.Range("M16").Value = 0.01485148514851485
$ws.Range("N16").Value = 0.004950495049504951
$ws.Range("O16").Value = 0.06930693069306931
$ws.Range("S16").Value = 0.1237623762376238
$ws.Range("F17").Value = 0.01642710472279261
$ws.Range("H17").Value = 0.1540041067761807
$ws.Range("I17").Value = 0.08213552361396304
$ws.Range("J17").Value = 0.431211498973306
$ws.Range("K17").Value = 0.1170431211498973
$ws.Range("M17").Value = 0.008213552361396304
$ws.Range("N17").Value = 0.002053388090349076
$ws.Range("O17").Value = 0.06776180698151951
$ws.Range("S17").Value = 0.1211498973305955
$ws.Range("F18").Value = 0.01612903225806452
$ws.Range("H18").Value = 0.1935483870967742
$ws.Range("I18").Value = 0.07526881720430108
$ws.Range("J18").Value = 0.4516129032258064
$ws.Range("K18").Value = 0.07526881720430108
$ws.Range("M18").Value = 0.01075268817204301
$ws.Range("O18").Value = 0.04838709677419355
$ws.Range("S18").Value = 0.1290322580645161
$ws.Range("F19").Value = 0.0238627889634601
$ws.Range("H19").Value = 0.2356450410141685
$ws.Range("I19").Value = 0.08873974645786727
$ws.Range("J19").Value = 0.3422818791946309
$ws.Range("K19").Value = 0.1029082774049217
$ws.Range("M19").Value = 0.02237136465324385
$ws.Range("O19").Value = 0.06487695749440715
$ws.Range("S19").Value = 0.1193139448173005
